# Manual Add student to project - Done
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Project 0")

$projectId = "5f8401d12db8b94ef525ec41"

# Row 4: Software role, 2 required, 0 left, new student enrolled
$ws.Range("A4").Value = $projectId
$ws.Range("B4").Value = "Software"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = '["5f7034d75f98895b90a4a029","testing",""]'

# Row 5: Engineer role, 2 required, 1 left, Mikaela enrolled
$ws.Range("A5").Value = $projectId
$ws.Range("B5").Value = "Engineer"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = '["Mikaela"]'

# Row 6: Data role, 2 required, 1 left, Mikaela enrolled
$ws.Range("A6").Value = $projectId
$ws.Range("B6").Value = "Data"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = '["Mikaela"]'
